# Delete rows that correspond to resampled / no-longer-present label files
# from each per-species training-output sheet.
#
# Sheet "s__Anaerotignum faecicola-b-p": rows 9-12 (1-based, incl. header)
#   i.e. labels UMGS137_11, UMGS137_18, UMGS137_2, UMGS137_6 dropped; the
#   remaining rows shift up, so the sheet ends at row 36 instead of 40.
#
# Sheet "s__Anaerotignum sp001304995-b-p": rows 9-18 (1-based, incl. header)
#   i.e. labels 12718_7_41_10, _14, _19, _20, _26, _27, _33, _35, _36, _4
#   dropped; remaining rows shift up, sheet ends at row 36 instead of 46.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("s__Anaerotignum faecicola-b-p")
$ws1.Range("A9:E12").EntireRow.Delete() | Out-Null

$ws2 = $wb.Worksheets.Item("s__Anaerotignum sp001304995-b-p")
$ws2.Range("A9:E18").EntireRow.Delete() | Out-Null
